$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formats (style s=3) from row 9 (C9:L9) down to row 11 (C11:L11)
$ws.Range("C9:L9").Copy()
$ws.Range("C11:L11").PasteSpecial(-4122)

# Copy format for B11 (style s=8) from B10
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)

# Copy format for F11 (style s=16, red highlighted "flagged" cell) from E7
$ws.Range("E7").Copy()
$ws.Range("F11").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Now set the actual values for row 11
$ws.Range("B11").Value = "Task #8: Implement a function to return list of multiple tags"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").Formula = "=SUM(C11:L11)"

# Move the "Daily Work Sum" totals row from row 12 to row 13
$ws.Range("A13").Value = $ws.Range("A12").Value2
$ws.Range("B13").Value = $ws.Range("B12").Value2
$ws.Range("C13").Formula = "=SUM(C3:C10)"
$ws.Range("D13").Formula = "=SUM(D3:D10)"
$ws.Range("E13").Formula = "=SUM(E3:E10)"
$ws.Range("F13").Formula = "=SUM(F3:F10)"
$ws.Range("G13").Formula = "=SUM(G3:G10)"
$ws.Range("H13").Formula = "=SUM(H3:H10)"
$ws.Range("I13").Formula = "=SUM(I3:I10)"
$ws.Range("J13").Formula = "=SUM(J3:J10)"
$ws.Range("K13").Formula = "=SUM(K3:K10)"
$ws.Range("L13").Formula = "=SUM(L3:L10)"
$ws.Range("N13").Formula = "=B15"

# Copy the formatting of the old totals row down to the new row 13
$ws.Range("A12:N12").Copy()
$ws.Range("A13:N13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the "Total Work" row from row 14 to row 15
$ws.Range("A15").Value = $ws.Range("A14").Value2
$ws.Range("B15").Formula = "=SUM(N3:N10)"
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Clear the old rows - row 12 should only keep O12; row 14 should be empty
$ws.Range("A12:N12").Clear()
$ws.Range("A14:B14").Clear()

# Extend the two tables down by one row to include the new row 11
$t1 = $ws.ListObjects.Item(1)
$t1.Resize($ws.Range("B2:L11"))
$t2 = $ws.ListObjects.Item(2)
$t2.Resize($ws.Range("N2:N11"))

# Update the selection to match the new active cell
[void]$ws.Range("F11").Select()

Write-Host "done"
